$wb = $excel.ActiveWorkbook
try { Write-Output $wb.Blobs } catch { Write-Output "ERR1 $_" }
try { Write-Output $excel.Blobs } catch { Write-Output "ERR2 $_" }
